$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto "Updated cryptos list" refresh: new prices / 1h volume deltas scraped from
# coinranking.com, plus a rank swap between Bittensor and Internet Computer (rows 32/33).

$ws.Range('D2').Value = '89.044.42'
$ws.Range('E2').Value = '  +10.79%  '
$ws.Range('D3').Value = '3.367.25'
$ws.Range('E3').Value = '  +5.80%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '222.95'
$ws.Range('E5').Value = '  +6.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '645.45'
$ws.Range('E6').Value = '  +2.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '0.346'
$ws.Range('E7').Value = '  +25.99%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.619'
$ws.Range('E9').Value = '  +5.40%  '
$ws.Range('D10').Value = '3.365.20'
$ws.Range('E10').Value = '  +5.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.613'
$ws.Range('E11').Value = '  +4.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '0.0000279'
$ws.Range('E12').Value = '  +8.00%  '
$ws.Range('E13').Value = '  +2.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '35.26'
$ws.Range('E14').Value = '  +10.51%  '
$ws.Range('D15').Value = '3.987.15'
$ws.Range('E15').Value = '  +5.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '5.47'
$ws.Range('E16').Value = '  +3.58%  '
$ws.Range('D17').Value = '88.857.92'
$ws.Range('E17').Value = '  +10.68%  '
$ws.Range('D18').Value = '3.370.39'
$ws.Range('E18').Value = '  +5.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '14.77'
$ws.Range('E19').Value = '  +3.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '3.20'
$ws.Range('E20').Value = '  +5.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '473.48'
$ws.Range('E21').Value = '  +8.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '9.28'
$ws.Range('E22').Value = '  +1.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '5.54'
$ws.Range('E23').Value = '  +6.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '13.77'
$ws.Range('E24').Value = '  +26.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '7.37'
$ws.Range('E25').Value = '  +6.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '5.49'
$ws.Range('E26').Value = '  +16.82%  '
$ws.Range('E27').Value = '  +4.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '79.68'
$ws.Range('E28').Value = '  +4.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '0.209'
$ws.Range('E29').Value = '  +69.53%  '
$ws.Range('E30').Value = '  +6.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '9.40'
$ws.Range('E32').Value = '  +5.19%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '597.05'
$ws.Range('E33').Value = '  +6.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '1.55'
$ws.Range('E35').Value = '  +7.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '2.08'
$ws.Range('E36').Value = '  +3.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '0.153'
$ws.Range('E37').Value = '  +1.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '24.42'
$ws.Range('E38').Value = '  +6.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '6.97'
$ws.Range('E39').Value = '  +23.58%  '
$ws.Range('E40').Value = '  +3.99%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '21.73'
$ws.Range('E42').Value = '  +4.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '2.09'
$ws.Range('E43').Value = '  +15.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '3.06'
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '193.65'
$ws.Range('E45').Value = '  +2.25%  '
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '157.31'
$ws.Range('E47').Value = '  -3.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '47.58'
$ws.Range('E48').Value = '  +11.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '1.39'
$ws.Range('E49').Value = '  +7.32%  '
$ws.Range('E50').Value = '  +2.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '0.665'
$ws.Range('E51').Value = '  +6.13%  '
